# Scheduled-runner refresh of market/profit figures on the Anima_Profits
# sheets (one per crafting class tab: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Updates cached currentAveragePrice* / LevePrice* / LeveProfit* columns
# (H:N) for the rows whose underlying market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 10767.4
$ws.Range("I19").Value = 758.3333
$ws.Range("J19").Value = 15057
$ws.Range("K19").Value = 758.3333
$ws.Range("L19").Value = 15057
$ws.Range("M19").Value = -583.3333
$ws.Range("N19").Value = -15407

$ws.Range("H33").Value = 494
$ws.Range("I33").Value = 494
$ws.Range("K33").Value = 494
$ws.Range("M33").Value = -265

$ws.Range("H40").Value = 2000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H64").Value = 3749.75
$ws.Range("J64").Value = 3333
$ws.Range("L64").Value = 3333
$ws.Range("N64").Value = -3829

$ws.Range("H67").Value = 3749.75
$ws.Range("J67").Value = 3333
$ws.Range("L67").Value = 3333
$ws.Range("N67").Value = -5049

$ws.Range("H86").Value = 114302024
$ws.Range("I86").Value = 133351780
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 133351780
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -133350657
$ws.Range("N86").Value = -5746

$ws.Range("H89").Value = 114302024
$ws.Range("I89").Value = 133351780
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 666758900
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -666753284
$ws.Range("N89").Value = -28732

$ws.Range("H92").Value = 16668302
$ws.Range("I92").Value = 20834502
$ws.Range("J92").Value = 3500
$ws.Range("K92").Value = 20834502
$ws.Range("L92").Value = 3500
$ws.Range("M92").Value = -20833254
$ws.Range("N92").Value = -5996

$ws.Range("H98").Value = 600.8125
$ws.Range("I98").Value = 607.5333000000001
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 607.5333000000001
$ws.Range("L98").Value = 500
$ws.Range("M98").Value = 890.4666999999999
$ws.Range("N98").Value = -3496

$ws.Range("H122").Value = 600.8125
$ws.Range("I122").Value = 607.5333000000001
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 1822.5999
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = 627.4000999999998
$ws.Range("N122").Value = -6400

$ws.Range("H138").Value = 2924.695
$ws.Range("I138").Value = 1559.6061
$ws.Range("J138").Value = 3844.0408
$ws.Range("K138").Value = 4678.8183
$ws.Range("L138").Value = 11532.1224
$ws.Range("M138").Value = 461.1817000000001
$ws.Range("N138").Value = -21812.1224

$ws.Range("H141").Value = 4219.2
$ws.Range("I141").Value = 1249.6666
$ws.Range("J141").Value = 8673.5
$ws.Range("K141").Value = 3748.9998
$ws.Range("L141").Value = 26020.5
$ws.Range("M141").Value = 1431.0002
$ws.Range("N141").Value = -36380.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 506326.5
$ws.Range("I32").Value = 552517.9399999999
$ws.Range("J32").Value = 24616.143
$ws.Range("K32").Value = 552517.9399999999
$ws.Range("L32").Value = 24616.143
$ws.Range("M32").Value = -552230.9399999999
$ws.Range("N32").Value = -25190.143

$ws.Range("H122").Value = 3951.0908
$ws.Range("I122").Value = 3180.2856
$ws.Range("J122").Value = 5300
$ws.Range("K122").Value = 9540.856800000001
$ws.Range("L122").Value = 15900
$ws.Range("M122").Value = -7090.856800000001
$ws.Range("N122").Value = -20800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 33335930
$ws.Range("I20").Value = 2594.9
$ws.Range("J20").Value = 100002600
$ws.Range("K20").Value = 2594.9
$ws.Range("L20").Value = 100002600
$ws.Range("M20").Value = -2347.9
$ws.Range("N20").Value = -100003094

$ws.Range("H134").Value = 2860.3684
$ws.Range("I134").Value = 2407.9524
$ws.Range("J134").Value = 3419.2354
$ws.Range("K134").Value = 7223.8572
$ws.Range("L134").Value = 10257.7062
$ws.Range("M134").Value = -4688.8572
$ws.Range("N134").Value = -15327.7062

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 172054.17
$ws.Range("I35").Value = 255881.25
$ws.Range("K35").Value = 255881.25
$ws.Range("M35").Value = -255587.25

$ws.Range("H58").Value = 1309.9429
$ws.Range("I58").Value = 1014
$ws.Range("J58").Value = 1753.8572
$ws.Range("K58").Value = 1014
$ws.Range("L58").Value = 1753.8572
$ws.Range("M58").Value = -811
$ws.Range("N58").Value = -2159.8572

$ws.Range("H62").Value = 3927.7778
$ws.Range("I62").Value = 3978.5715
$ws.Range("J62").Value = 3750
$ws.Range("K62").Value = 3978.5715
$ws.Range("L62").Value = 3750
$ws.Range("M62").Value = -3354.5715
$ws.Range("N62").Value = -4998

$ws.Range("H65").Value = 3927.7778
$ws.Range("I65").Value = 3978.5715
$ws.Range("J65").Value = 3750
$ws.Range("K65").Value = 19892.8575
$ws.Range("L65").Value = 18750
$ws.Range("M65").Value = -16772.8575
$ws.Range("N65").Value = -24990

$ws.Range("H132").Value = 35355252
$ws.Range("I132").Value = 45456100
$ws.Range("J132").Value = 15153556
$ws.Range("K132").Value = 136368300
$ws.Range("L132").Value = 45460668
$ws.Range("M132").Value = -136365770
$ws.Range("N132").Value = -45465728

$ws.Range("H134").Value = 5152.625
$ws.Range("I134").Value = 5120.0386
$ws.Range("J134").Value = 5293.8335
$ws.Range("K134").Value = 15360.1158
$ws.Range("L134").Value = 15881.5005
$ws.Range("M134").Value = -12825.1158
$ws.Range("N134").Value = -20951.5005

$ws.Range("H136").Value = 1309.9429
$ws.Range("I136").Value = 1014
$ws.Range("J136").Value = 1753.8572
$ws.Range("K136").Value = 3042
$ws.Range("L136").Value = 5261.571599999999
$ws.Range("M136").Value = -492
$ws.Range("N136").Value = -10361.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 238.2
$ws.Range("I47").Value = 238.2
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 714.5999999999999
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -283.5999999999999

$ws.Range("H87").Value = 3195.4
$ws.Range("I87").Value = 725.6667
$ws.Range("K87").Value = 2177.0001
$ws.Range("M87").Value = -929.0001000000002

$ws.Range("H90").Value = 3195.4
$ws.Range("I90").Value = 725.6667
$ws.Range("K90").Value = 6531.0003
$ws.Range("M90").Value = -291.0002999999997

$ws.Range("H107").Value = 43478604
$ws.Range("I107").Value = 389
$ws.Range("J107").Value = 76923384
$ws.Range("K107").Value = 1167
$ws.Range("L107").Value = 230770152
$ws.Range("M107").Value = 753
$ws.Range("N107").Value = -230773992

$ws.Range("H120").Value = 11861.429
$ws.Range("I120").Value = 11515
$ws.Range("K120").Value = 34545
$ws.Range("M120").Value = -29707

$ws.Range("H122").Value = 3149.3845
$ws.Range("I122").Value = 331.54544
$ws.Range("J122").Value = 6796
$ws.Range("K122").Value = 2983.90896
$ws.Range("L122").Value = 61164
$ws.Range("M122").Value = -533.9089599999998
$ws.Range("N122").Value = -66064

$ws.Range("H131").Value = 7717.3335
$ws.Range("I131").Value = 503.75
$ws.Range("J131").Value = 9778.357
$ws.Range("K131").Value = 1511.25
$ws.Range("L131").Value = 29335.071
$ws.Range("M131").Value = 3528.75
$ws.Range("N131").Value = -39415.071

$ws.Range("H138").Value = 3598.3333
$ws.Range("I138").Value = 1338.3846
$ws.Range("J138").Value = 7270.75
$ws.Range("K138").Value = 4015.1538
$ws.Range("L138").Value = 21812.25
$ws.Range("M138").Value = 1124.8462
$ws.Range("N138").Value = -32092.25

$ws.Range("H139").Value = 3629.8235
$ws.Range("I139").Value = 1654.6666
$ws.Range("J139").Value = 5189.1577
$ws.Range("K139").Value = 4963.9998
$ws.Range("L139").Value = 15567.4731
$ws.Range("M139").Value = 176.0002000000004
$ws.Range("N139").Value = -25847.4731

$ws.Range("H141").Value = 7405.727
$ws.Range("I141").Value = 2365.75
$ws.Range("J141").Value = 10285.714
$ws.Range("K141").Value = 7097.25
$ws.Range("L141").Value = 30857.142
$ws.Range("M141").Value = -1917.25
$ws.Range("N141").Value = -41217.142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8925.091
$ws.Range("I70").Value = 10340.8125
$ws.Range("J70").Value = 5149.8335
$ws.Range("K70").Value = 10340.8125
$ws.Range("L70").Value = 5149.8335
$ws.Range("M70").Value = -10070.8125
$ws.Range("N70").Value = -5689.8335

$ws.Range("H73").Value = 8925.091
$ws.Range("I73").Value = 10340.8125
$ws.Range("J73").Value = 5149.8335
$ws.Range("K73").Value = 10340.8125
$ws.Range("L73").Value = 5149.8335
$ws.Range("M73").Value = -9404.8125
$ws.Range("N73").Value = -7021.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4285.7144
$ws.Range("I46").Value = 1700
$ws.Range("J46").Value = 10750
$ws.Range("K46").Value = 1700
$ws.Range("L46").Value = 10750
$ws.Range("M46").Value = -1512
$ws.Range("N46").Value = -11126

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120
